# Updates the stack-trace text in the "setNumbering invalid level" expected
# generation document so it reflects the refreshed library/call-stack line
# numbers that come with moving the test dependencies from 3.1.0 to 3.1.1
# (object hash in the MParagraphImpl message also changes since it is a
# fresh JVM object identity each run).
$d = $word.ActiveDocument

$oldStackTrace = "setNumbering(org.obeonetwork.m2doc.element.MParagraph,java.lang.Integer,java.lang.Integer) with arguments [org.obeonetwork.m2doc.element.impl.MParagraphImpl@20fa5277, 1, 9] failed:`n`tmaximum level for numbering with ID 1 is 8, 9 is too high.`njava.lang.IllegalArgumentException: maximum level for numbering with ID 1 is 8, 9 is too high.`n`tat org.obeonetwork.m2doc.services.PaginationServices.setNumbering(PaginationServices.java:410)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:281)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:805)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)`n`tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n"
$newStackTrace = "setNumbering(org.obeonetwork.m2doc.element.MParagraph,java.lang.Integer,java.lang.Integer) with arguments [org.obeonetwork.m2doc.element.impl.MParagraphImpl@66abb2fa, 1, 9] failed:`n`tmaximum level for numbering with ID 1 is 8, 9 is too high.`njava.lang.IllegalArgumentException: maximum level for numbering with ID 1 is 8, 9 is too high.`n`tat org.obeonetwork.m2doc.services.PaginationServices.setNumbering(PaginationServices.java:410)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:587)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:282)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:845)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)`n`tat sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n"

$find = $d.Content.Find
$result = $find.Execute($oldStackTrace, $true, $false, $false, $false, $false, $true, 1, $false, $newStackTrace, 2)

if (-not $result) {
    throw "Could not find the expected stack trace text to replace."
}
